$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (so it lands at the end)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "nr_studies"

# Header row (bold + centered, matching the other sheets' header style)
$ws.Range("A1").Value = "outcome"
$ws.Range("B1").Value = "country_id_europe_1_usa_north_america_2_asia_3_australia_4_south_america_5"
$ws.Range("C1").Value = "n_effect_sizes"
$ws.Range("D1").Value = "k_studies"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108

function Set-Row($row, $outcome, $countryId, $nEffectSizes, $kStudies) {
    $ws.Range("A$row").Value = $outcome
    if ($countryId -ne $null) {
        $ws.Range("B$row").NumberFormat = "@"
        $ws.Range("B$row").Value = $countryId
    }
    $ws.Range("C$row").Value = $nEffectSizes
    $ws.Range("D$row").Value = $kStudies
}

Set-Row 2  "NS" "1"  460 48
Set-Row 3  "NS" "4"  19  2
Set-Row 4  "NS" $null 56 6
Set-Row 5  "NS" "3"  43  11
Set-Row 6  "NS" "2"  139 12
Set-Row 7  "NS" "5"  5   2
Set-Row 8  "NT" "1"  253 35
Set-Row 9  "NT" "3"  47  8
Set-Row 10 "NT" "2"  51  3
Set-Row 11 "NT" $null 15 2
Set-Row 12 "NT" "5"  5   2
Set-Row 13 "NT" "4"  10  1
